$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.086.79"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "3.806.20"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +12.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.07%  "
$ws.Range("D7").Value = "3.806.25"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.87%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +7.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.81%  "
$ws.Range("D15").Value = "4.444.44"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.802.67"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "71.109.12"
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +18.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.53%  "
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("D29").Value = "3.956.50"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +17.35%  "
$ws.Range("E32").Value = "  +8.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.32"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.88%  "
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "3.756.80"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("E40").Value = "  +6.77%  "
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("E42").Value = "  +14.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000331"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +25.46%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.74"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "161.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "49.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.39%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.302"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.04%  "
